# Author deleted a single blog entry ("「我が魂と共に」مع نفسي ...") that
# occupied row 858 of Sheet1. Deleting the whole row shifts every
# subsequent row (859-883) up by one so the sheet now ends at row 882,
# which also updates the sheet's used-range/dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(858).Delete()
